# Atualizado por script em 01-12-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 42
$prev = $row - 1

# Copy formatting (style) from the row above, matching column by column
$ws.Range("A$prev`:V$prev").Copy()
$ws.Range("A$row`:V$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 41
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45261.64583333334
$ws.Cells.Item($row, 6).Value = "Jamshedpur"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Odisha FC"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 2.63
$ws.Cells.Item($row, 11).Value = "27/11/2023 15:42"
$ws.Cells.Item($row, 12).Value = 2.56
$ws.Cells.Item($row, 13).Value = "01/12/2023 15:22"
$ws.Cells.Item($row, 14).Value = 3.34
$ws.Cells.Item($row, 15).Value = "27/11/2023 15:42"
$ws.Cells.Item($row, 16).Value = 3.37
$ws.Cells.Item($row, 17).Value = "01/12/2023 15:29"
$ws.Cells.Item($row, 18).Value = 2.67
$ws.Cells.Item($row, 19).Value = "27/11/2023 15:42"
$ws.Cells.Item($row, 20).Value = 2.81
$ws.Cells.Item($row, 21).Value = "01/12/2023 15:22"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/jamshedpur-odisha-fc/bL9qkTaF/"
